$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the "用户设置" (user settings) sheet: insert the new
#    "api-requestCode" step before the email-check block, rename
#    "核对email" -> "核对旧email", and append a new "手机号 Tab"
#    (phone) block before the final "api - updatePhone" row.
# ------------------------------------------------------------------
$settings = $wb.Worksheets.Item(3)

# Make room for the new "api-requestCode" row (old row 22 -> 24).
$settings.Rows("22:23").Insert() | Out-Null
$settings.Range("D22").Value = "api-requestCode"

# Make room for the new phone block before the old row 31 (now at 33).
$settings.Rows("31").Insert() | Out-Null

# "核对email" (now at row 28) becomes "核对旧email".
$settings.Range("G28").Value = "核对旧email"

# New "手机号 Tab" block.
$settings.Range("C32").Value = "手机号 Tab"
$settings.Range("D33").Value = "更改手机号，需更新Token，因为，Token中Claim 手机号"
$settings.Range("D33").Font.Color = 255

$settings.Range("E27").Select() | Out-Null

# The "登录" (login) sheet's last active cell also moved.
$login = $wb.Worksheets.Item(2)
$login.Range("D12").Select() | Out-Null

# ------------------------------------------------------------------
# 2) Add a new "Store" sheet as the first tab, describing the Vuex
#    store shape used by the app.
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$store = $wb.Worksheets.Add($firstSheet)
$store.Name = "Store"

$store.Range("B3").Value = "login"
$store.Range("B1").Value = "setter"
$store.Range("G1").Value = "clear"
$store.Range("G3").Value = "logout"
$store.Range("C4").Value = "store {active, detailed_info}"

$store.Activate() | Out-Null
$store.Range("G3").Select() | Out-Null
